$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, pushing the existing rows 38-56 down to 39-57
$ws.Rows.Item(38).Insert()

# Populate the newly-inserted row 38 with the new weekly data record
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(38, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value = 44992
$ws.Cells.Item(38, 5).Value = 15
$ws.Cells.Item(38, 6).Value = 100112045
$ws.Cells.Item(38, 7).Value = "Zapallo"
$ws.Cells.Item(38, 8).Value = "Camote"
$ws.Cells.Item(38, 9).Value = "1a (cosecha)"
$ws.Cells.Item(38, 10).Value = 800
$ws.Cells.Item(38, 11).Value = 630
$ws.Cells.Item(38, 12).Value = 650
$ws.Cells.Item(38, 13).Value = 640
$ws.Cells.Item(38, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(38, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(38, 16).Value = 640
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"
